# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Corrige los valores de Periodo Mora y Valor Mora de los trabajadores en
# las filas 17 y 18 (se habian intercambiado los periodos 1905/1906).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fila 17 (MARINA ESTER JIMENEZ ZABALA, periodo 1906 -> 1905)
$ws.Range("E17").Value = "1905"
$ws.Range("F17").Value = 27667

# Fila 18 (MARINA ESTER JIMENEZ ZABALA, periodo 1905 -> 1906)
$ws.Range("E18").Value = "1906"
$ws.Range("F18").Value = 6640
